# Added nonmilitary-industry data: the sparse year rows (2-18) get filled out
# with the full set of BLS series (D:M), and the leftover scattered rows
# 19-33 -- which held overflow of that same data under the old sparse
# layout -- are cleared now that everything lives in one dense row per year.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C7 had a stray value (883, a fat-fingered duplicate of the women's figure)
# that should read 88.3
$ws.Range("C7").Value = 88.3

# Fill in D:M for rows 2-18 per target data (dense consolidated table)
# Row 2
$ws.Range("D2").Value = 515.29999999999995
$ws.Range("E2").Value = 101.25

# Row 3
$ws.Range("D3").Value = 476.2
$ws.Range("E3").Value = 106.63
$ws.Range("F3").Value = 2.62
$ws.Range("G3").Value = 40.700000000000003
$ws.Range("H3").Value = 2.1
$ws.Range("I3").Value = 2.2999999999999998
$ws.Range("K3").Value = 3.2
$ws.Range("L3").Value = 1.2
$ws.Range("M3").Value = 1.5

# Row 4
$ws.Range("D4").Value = 407.4
$ws.Range("E4").Value = 110.43
$ws.Range("F4").Value = 2.7
$ws.Range("H4").Value = 2.2000000000000002
$ws.Range("I4").Value = 2.4
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.7

# Row 5
$ws.Range("D5").Value = 351.7
$ws.Range("E5").Value = 114.68
$ws.Range("F5").Value = 2.77
$ws.Range("G5").Value = 41.4
$ws.Range("H5").Value = 2.5
$ws.Range("I5").Value = 2.8
$ws.Range("J5").Value = 1.7
$ws.Range("K5").Value = 2.6
$ws.Range("L5").Value = 0.9
$ws.Range("M5").Value = 1.3

# Row 6
$ws.Range("D6").Value = 356.4
$ws.Range("E6").Value = 119.97
$ws.Range("F6").Value = 2.87
$ws.Range("G6").Value = 41.8
$ws.Range("H6").Value = 2.9
$ws.Range("I6").Value = 2.9
$ws.Range("J6").Value = 2.1
$ws.Range("K6").Value = 2.5
$ws.Range("L6").Value = 1.1000000000000001
$ws.Range("M6").Value = 1

# Row 7
$ws.Range("D7").Value = 358.2
$ws.Range("E7").Value = 122.43
$ws.Range("F7").Value = 2.95
$ws.Range("G7").Value = 41.5
$ws.Range("H7").Value = 2.6
$ws.Range("I7").Value = 2.4
$ws.Range("J7").Value = 1.7
$ws.Range("K7").Value = 2.5
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.1000000000000001

# Row 8
$ws.Range("D8").Value = 352.7
$ws.Range("E8").Value = 125.03
$ws.Range("F8").Value = 3.02
$ws.Range("G8").Value = 41.4
$ws.Range("H8").Value = 2.5
$ws.Range("I8").Value = 2.1
$ws.Range("J8").Value = 1.4
$ws.Range("K8").Value = 2.6
$ws.Range("L8").Value = 0.9
$ws.Range("M8").Value = 1.2

# Row 9
$ws.Range("D9").Value = 339.2
$ws.Range("E9").Value = 131.88
$ws.Range("F9").Value = 3.14
$ws.Range("G9").Value = 42
$ws.Range("H9").Value = 3.3
$ws.Range("I9").Value = 3.2
$ws.Range("J9").Value = 2.5
$ws.Range("K9").Value = 2.2999999999999998
$ws.Range("L9").Value = 1.1000000000000001
$ws.Range("M9").Value = 0.7

# Row 10
$ws.Range("D10").Value = 402.5
$ws.Range("E10").Value = 143.32
$ws.Range("F10").Value = 3.31
$ws.Range("G10").Value = 43.3
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 4.2
$ws.Range("J10").Value = 3.6
$ws.Range("K10").Value = 2.6
$ws.Range("L10").Value = 1.6
$ws.Range("M10").Value = 0.3

# Row 11
$ws.Range("D11").Value = 489.5
$ws.Range("E11").Value = 146.97
$ws.Range("F11").Value = 3.45
$ws.Range("G11").Value = 42.6
$ws.Range("H11").Value = 4.4000000000000004
$ws.Range("I11").Value = 3.1
$ws.Range("J11").Value = 2.6
$ws.Range("K11").Value = 2.7
$ws.Range("L11").Value = 1.6
$ws.Range("M11").Value = 5

# Row 12
$ws.Range("D12").Value = 519.79999999999995
$ws.Range("E12").Value = 152.04
$ws.Range("F12").Value = 3.62
$ws.Range("G12").Value = 42
$ws.Range("H12").Value = 3.8
$ws.Range("I12").Value = 2.2999999999999998
$ws.Range("J12").Value = 1.8
$ws.Range("K12").Value = 2.7
$ws.Range("L12").Value = 1.5
$ws.Range("M12").Value = 0.7

# Row 13
$ws.Range("D13").Value = 485.6
$ws.Range("E13").Value = 161.35
$ws.Range("F13").Value = 3.86
$ws.Range("G13").Value = 41.8
$ws.Range("H13").Value = 3.4
$ws.Range("I13").Value = 2
$ws.Range("J13").Value = 1.4
$ws.Range("K13").Value = 2.8
$ws.Range("L13").Value = 1.3
$ws.Range("M13").Value = 1

# Row 14
$ws.Range("D14").Value = 421.9
$ws.Range("E14").Value = 168.51
$ws.Range("F14").Value = 4.1100000000000003
$ws.Range("G14").Value = 41
$ws.Range("H14").Value = 2.7
$ws.Range("I14").Value = 1.4
$ws.Range("J14").Value = 0.7
$ws.Range("K14").Value = 3.5
$ws.Range("L14").Value = 0.8
$ws.Range("M14").Value = 2.2000000000000002

# Row 15
$ws.Range("D15").Value = 321.7
$ws.Range("E15").Value = 175.82
$ws.Range("F15").Value = 4.32
$ws.Range("G15").Value = 40.700000000000003
$ws.Range("H15").Value = 2.2999999999999998
$ws.Range("I15").Value = 1.7
$ws.Range("J15").Value = 0.7
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 0.6
$ws.Range("M15").Value = 1.9

# Row 16
$ws.Range("D16").Value = 270
$ws.Range("E16").Value = 193.44
$ws.Range("F16").Value = 4.6500000000000004
$ws.Range("G16").Value = 41.6
$ws.Range("H16").Value = 3
$ws.Range("I16").Value = 2
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.7
$ws.Range("M16").Value = 9

# Row 17
$ws.Range("D17").Value = 287.3
$ws.Range("E17").Value = 207.92
$ws.Range("F17").Value = 5.01
$ws.Range("G17").Value = 41.5
$ws.Range("H17").Value = 3.4
$ws.Range("I17").Value = 2.2999999999999998
$ws.Range("J17").Value = 1.5
$ws.Range("K17").Value = 2.1
$ws.Range("L17").Value = 0.9
$ws.Range("M17").Value = 0.7

# Row 18
$ws.Range("D18").Value = 295.10000000000002
$ws.Range("E18").Value = 218.7
$ws.Range("F18").Value = 5.4
$ws.Range("G18").Value = 40.5
$ws.Range("H18").Value = 3.2
$ws.Range("I18").Value = 2.1
$ws.Range("J18").Value = 1.5
$ws.Range("K18").Value = 1.9
$ws.Range("L18").Value = 0.8
$ws.Range("M18").Value = 0.5

# Remove old rows 19-33 (their data has been consolidated into rows 2-18 above)
$ws.Range("A19:M33").ClearContents()

# Update selection to C7 (matches target sheetView state)
$ws.Range("C7").Select()